$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 223
$ws.Range("I6").Value = 196.07143
$ws.Range("K6").Value = 588.21429
$ws.Range("M6").Value = -476.21429
# Row 135
$ws.Range("H135").Value = 8621748
$ws.Range("I135").Value = 10417528
$ws.Range("J135").Value = 2007.2
$ws.Range("K135").Value = 93757752
$ws.Range("L135").Value = 18064.8
$ws.Range("M135").Value = -93755217
$ws.Range("N135").Value = -23134.8
# Row 139
$ws.Range("H139").Value = 44500
$ws.Range("J139").Value = 44500
$ws.Range("L139").Value = 44500
$ws.Range("N139").Value = -54780

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 3
$ws.Range("H3").Value = 2236.111
$ws.Range("J3").Value = 2842.8572
$ws.Range("L3").Value = 2842.8572
$ws.Range("N3").Value = -3072.8572
# Row 61
$ws.Range("H61").Value = 13890349
$ws.Range("I61").Value = 13890349
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 13890349
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -13890137
$ws.Range("N61").ClearContents() | Out-Null
# Row 74
$ws.Range("H74").Value = 1645.76
$ws.Range("J74").Value = 1449.3077
$ws.Range("L74").Value = 1449.3077
$ws.Range("N74").Value = -3197.3077
# Row 77
$ws.Range("H77").Value = 1645.76
$ws.Range("J77").Value = 1449.3077
$ws.Range("L77").Value = 7246.538500000001
$ws.Range("N77").Value = -15982.5385
# Row 130
$ws.Range("H130").Value = 30429
$ws.Range("J130").Value = 30429
$ws.Range("L130").Value = 30429
$ws.Range("N130").Value = -40469
# Row 132
$ws.Range("H132").Value = 2101888.5
$ws.Range("I132").Value = 980.3158
$ws.Range("J132").Value = 6537139.5
$ws.Range("K132").Value = 2940.9474
$ws.Range("L132").Value = 19611418.5
$ws.Range("M132").Value = -410.9474
$ws.Range("N132").Value = -19616478.5
# Row 136
$ws.Range("H136").Value = 13890349
$ws.Range("I136").Value = 13890349
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 41671047
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -41668497
$ws.Range("N136").ClearContents() | Out-Null

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 5
$ws.Range("H5").Value = 3018
$ws.Range("I5").Value = 29.666666
$ws.Range("J5").Value = 4811
$ws.Range("K5").Value = 29.666666
$ws.Range("L5").Value = 4811
$ws.Range("M5").Value = 83.33333400000001
$ws.Range("N5").Value = -5037
# Row 105
$ws.Range("H105").Value = 83334620
$ws.Range("I105").Value = 1410.1818
$ws.Range("K105").Value = 1410.1818
$ws.Range("M105").Value = 336.8181999999999
# Row 128
$ws.Range("H128").Value = 2454
$ws.Range("I128").Value = 2454
$ws.Range("K128").Value = 7362
$ws.Range("M128").Value = -4872

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 4149.8423
$ws.Range("I3").Value = 3256.4666
$ws.Range("J3").Value = 7500
$ws.Range("K3").Value = 9769.399800000001
$ws.Range("L3").Value = 22500
$ws.Range("M3").Value = -9657.399800000001
$ws.Range("N3").Value = -22724
# Row 97
$ws.Range("H97").Value = 1683.3334
$ws.Range("I97").Value = 800
$ws.Range("K97").Value = 2400
$ws.Range("M97").Value = -1904
# Row 113
$ws.Range("H113").Value = 11945155
$ws.Range("I113").Value = 5556361.5
$ws.Range("J113").Value = 31111534
$ws.Range("K113").Value = 16669084.5
$ws.Range("L113").Value = 93334602
$ws.Range("M113").Value = -16666914.5
$ws.Range("N113").Value = -93338942
# Row 115
$ws.Range("H115").Value = 4975.5713
$ws.Range("I115").Value = 500.4
$ws.Range("J115").Value = 6374.0625
$ws.Range("K115").Value = 1501.2
$ws.Range("L115").Value = 19122.1875
$ws.Range("M115").Value = -326.1999999999998
$ws.Range("N115").Value = -21472.1875
# Row 131
$ws.Range("H131").Value = 782.83
$ws.Range("J131").Value = 782.83
$ws.Range("L131").Value = 2348.49
$ws.Range("N131").Value = -12428.49
# Row 138
$ws.Range("H138").Value = 17545590
$ws.Range("I138").Value = 25641902
$ws.Range("J138").Value = 3583.1667
$ws.Range("K138").Value = 76925706
$ws.Range("L138").Value = 10749.5001
$ws.Range("M138").Value = -76920566
$ws.Range("N138").Value = -21029.5001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 1805.409
$ws.Range("I102").Value = 909.25
$ws.Range("J102").Value = 2880.8
$ws.Range("K102").Value = 909.25
$ws.Range("L102").Value = 2880.8
$ws.Range("M102").Value = 712.75
$ws.Range("N102").Value = -6124.8
# Row 122
$ws.Range("H122").Value = 31257838
$ws.Range("I122").Value = 38470624
$ws.Range("J122").Value = 2433.3333
$ws.Range("K122").Value = 115411872
$ws.Range("L122").Value = 7299.999899999999
$ws.Range("M122").Value = -115409422
$ws.Range("N122").Value = -12199.9999
# Row 132
$ws.Range("H132").Value = 20101.666
$ws.Range("I132").Value = 2806
$ws.Range("K132").Value = 8418
$ws.Range("M132").Value = -5888

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 1888.9445
$ws.Range("I7").Value = 1181.7273
$ws.Range("J7").Value = 3000.2856
$ws.Range("K7").Value = 1181.7273
$ws.Range("L7").Value = 3000.2856
$ws.Range("M7").Value = -1069.7273
$ws.Range("N7").Value = -3224.2856
# Row 26
$ws.Range("H26").Value = 850
$ws.Range("J26").Value = 850
$ws.Range("L26").Value = 850
$ws.Range("N26").Value = -1440
# Row 62
$ws.Range("H62").Value = 5225.5
$ws.Range("I62").Value = 5225.5
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 5225.5
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -4601.5
$ws.Range("N62").ClearContents() | Out-Null
# Row 65
$ws.Range("H65").Value = 5225.5
$ws.Range("I65").Value = 5225.5
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 15676.5
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -12556.5
$ws.Range("N65").ClearContents() | Out-Null
# Row 126
$ws.Range("H126").Value = 1888.9445
$ws.Range("I126").Value = 1181.7273
$ws.Range("J126").Value = 3000.2856
$ws.Range("K126").Value = 3545.1819
$ws.Range("L126").Value = 9000.856800000001
$ws.Range("M126").Value = -1075.1819
$ws.Range("N126").Value = -13940.8568
# Row 132
$ws.Range("H132").Value = 114310510
$ws.Range("J132").Value = 30999.625
$ws.Range("L132").Value = 92998.875
$ws.Range("N132").Value = -98058.875

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 820.2
$ws.Range("I107").Value = 500
$ws.Range("J107").Value = 900.25
$ws.Range("K107").Value = 1500
$ws.Range("L107").Value = 2700.75
$ws.Range("M107").Value = 420
$ws.Range("N107").Value = -6540.75
# Row 126
$ws.Range("H126").Value = 1453.4
$ws.Range("I126").Value = 748.75
$ws.Range("J126").Value = 2510.375
$ws.Range("K126").Value = 2246.25
$ws.Range("L126").Value = 7531.125
$ws.Range("M126").Value = 223.75
$ws.Range("N126").Value = -12471.125
# Row 132
$ws.Range("H132").Value = 35308.637
$ws.Range("I132").Value = 114056.445
$ws.Range("J132").Value = 5778.2085
$ws.Range("K132").Value = 342169.335
$ws.Range("L132").Value = 17334.6255
$ws.Range("M132").Value = -339639.335
$ws.Range("N132").Value = -22394.6255
